$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.208529472351074
$ws.Range("B1").Value = 2.64684009552002
$ws.Range("C1").Value = 1.907177925109863
$ws.Range("D1").Value = 1.758663177490234
$ws.Range("E1").Value = 1.811660170555115
